# Update crypto price/volume figures per the latest symbol-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.43"
$ws.Range("E2").Value = "'1.85%"
$ws.Range("D3").Value = "'36.01"
$ws.Range("E3").Value = "'2.92%"
$ws.Range("D4").Value = "'5.102"
$ws.Range("E4").Value = "'1.75%"
$ws.Range("D5").Value = "'0.08110"
$ws.Range("E5").Value = "'2.67%"
$ws.Range("D6").Value = "'1.953"
$ws.Range("E6").Value = "'1.18%"
$ws.Range("D7").Value = "'7.762"
$ws.Range("E7").Value = "'0.19%"
$ws.Range("D8").Value = "'0.9327"
$ws.Range("E8").Value = "'1.13%"
$ws.Range("D9").Value = "'0.1419"
$ws.Range("E9").Value = "'21.62%"
$ws.Range("D10").Value = "'0.1918"
$ws.Range("E10").Value = "'4.56%"
$ws.Range("D11").Value = "'0.09234"
$ws.Range("E11").Value = "'-1.12%"
$ws.Range("D12").Value = "'0.03530"
$ws.Range("E12").Value = "'-0.18%"
$ws.Range("D13").Value = "'0.09863"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("D14").Value = "'0.001409"
$ws.Range("E14").Value = "'0.82%"
$ws.Range("D15").Value = "'0.005874"
$ws.Range("E15").Value = "'1.20%"
$ws.Range("D16").Value = "'3.593"
$ws.Range("E16").Value = "'2.75%"
$ws.Range("D17").Value = "'4.192"
$ws.Range("E17").Value = "'4.17%"
$ws.Range("D18").Value = "'2.986"
$ws.Range("E18").Value = "'0.73%"
$ws.Range("D19").Value = "'0.3441"
$ws.Range("E19").Value = "'-0.06%"
$ws.Range("D20").Value = "'0.1350"
$ws.Range("E20").Value = "'3.20%"
$ws.Range("D21").Value = "'4.888"
$ws.Range("D22").Value = "'0.2412"
$ws.Range("E22").Value = "'0.58%"
$ws.Range("D23").Value = "'0.04513"
$ws.Range("E23").Value = "'0.37%"
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'0.32%"
$ws.Range("D25").Value = "'0.004868"
$ws.Range("E25").Value = "'6.49%"
$ws.Range("D26").Value = "'0.0001240"
$ws.Range("E26").Value = "'-0.81%"
$ws.Range("D39").Value = "'0.02003"
$ws.Range("E39").Value = "'5.82%"
$ws.Range("D40").Value = "'0.04941"
$ws.Range("E40").Value = "'5.27%"
$ws.Range("D41").Value = "'0.01078"
$ws.Range("E41").Value = "'12.98%"
$ws.Range("D42").Value = "'0.007648"
$ws.Range("E42").Value = "'0.72%"
$ws.Range("D43").Value = "'0.1381"
$ws.Range("E43").Value = "'4.44%"
$ws.Range("D44").Value = "'0.002100"
$ws.Range("E44").Value = "'-0.96%"
$ws.Range("D45").Value = "'0.01004"
$ws.Range("E45").Value = "'-10.08%"
$ws.Range("D46").Value = "'0.00006456"
$ws.Range("E46").Value = "'7.41%"
$ws.Range("E47").Value = "'0.22%"
$ws.Range("D49").Value = "'0.001193"
$ws.Range("E49").Value = "'-8.56%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.22%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.22%"
